$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================
# Step 1: stamp per-row formatting from template rows BEFORE
# writing any values (PasteSpecial -4122 = xlPasteFormats, so
# only styles move, no cell contents).
#   row parity 4/5/6  <- template row 2   (A2:L2 / M2 / N2)
#   row parity 7/8/9  <- template row 3   (A3:L3 / M3 / N3)
#   final new row     <- template row 143 (A143:L143, style 15/16)
#                        + M271/N271-sourced style 19 for the two
#                        trailing blank cells M283/N283
# ============================================================

# the last row (283) reuses the "blank trailing cell" style (19)
# that M271 currently has -- grab it first, before M271 is cleared
$ws.Range("M271").Copy()
$ws.Range("M283").PasteSpecial(-4122)
$ws.Range("M271").Copy()
$ws.Range("N283").PasteSpecial(-4122)

$ws.Range("A2:L2").Copy()
$ws.Range("A272:L272").PasteSpecial(-4122)
$ws.Range("A3:L3").Copy()
$ws.Range("A273:L273").PasteSpecial(-4122)
$ws.Range("M3").Copy()
$ws.Range("M273").PasteSpecial(-4122)
$ws.Range("A2:L2").Copy()
$ws.Range("A274:L274").PasteSpecial(-4122)
$ws.Range("A3:L3").Copy()
$ws.Range("A275:L275").PasteSpecial(-4122)
$ws.Range("M3").Copy()
$ws.Range("M275").PasteSpecial(-4122)
$ws.Range("A2:L2").Copy()
$ws.Range("A276:L276").PasteSpecial(-4122)
$ws.Range("N2").Copy()
$ws.Range("N276").PasteSpecial(-4122)
$ws.Range("A3:L3").Copy()
$ws.Range("A277:L277").PasteSpecial(-4122)
$ws.Range("M3").Copy()
$ws.Range("M277").PasteSpecial(-4122)
$ws.Range("A2:L2").Copy()
$ws.Range("A278:L278").PasteSpecial(-4122)
$ws.Range("N2").Copy()
$ws.Range("N278").PasteSpecial(-4122)
$ws.Range("A3:L3").Copy()
$ws.Range("A279:L279").PasteSpecial(-4122)
$ws.Range("A2:L2").Copy()
$ws.Range("A280:L280").PasteSpecial(-4122)
$ws.Range("M2").Copy()
$ws.Range("M280").PasteSpecial(-4122)
$ws.Range("A3:L3").Copy()
$ws.Range("A281:L281").PasteSpecial(-4122)
$ws.Range("N3").Copy()
$ws.Range("N281").PasteSpecial(-4122)
$ws.Range("A2:L2").Copy()
$ws.Range("A282:L282").PasteSpecial(-4122)
$ws.Range("A143:L143").Copy()
$ws.Range("A283:L283").PasteSpecial(-4122)

# ============================================================
# Step 2: set the actual cell values for the new rows
# ============================================================
$ws.Range("A272").Value = 45599.434015277773
$ws.Range("B272").Value = 'wnsrl2498@naver.com'
$ws.Range("C272").Value = '체육학과'
$ws.Range("D272").Value = 20227093
$ws.Range("E272").Value = '김준기'
$ws.Range("F272").Value = '스페인'
$ws.Range("G272").Value = '취업자 / 15세 이상 인구'
$ws.Range("H272").Value = '조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자'
$ws.Range("I272").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("J272").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K272").Value = '"15%"'
$ws.Range("L272").Value = 'Red'
$ws.Range("A273").Value = 45599.442972465273
$ws.Range("B273").Value = 'sihyune1104@gmail.com'
$ws.Range("C273").Value = '법학과'
$ws.Range("D273").Value = 20232726
$ws.Range("E273").Value = '안시현'
$ws.Range("F273").Value = '대한민국'
$ws.Range("G273").Value = '취업자 / 15세 이상 인구'
$ws.Range("H273").Value = '조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자'
$ws.Range("I273").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J273").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K273").Value = '"19.1%"'
$ws.Range("L273").Value = 'Red'
$ws.Range("M273").Value = '나는 사후 장기기증에 참여할 뜻이 없다'
$ws.Range("A274").Value = 45599.503719618056
$ws.Range("B274").Value = 'leedowon567@naver.com'
$ws.Range("C274").Value = '일본학과'
$ws.Range("D274").Value = 20211625
$ws.Range("E274").Value = '이도원'
$ws.Range("F274").Value = '대한민국'
$ws.Range("G274").Value = '경제활동인구 / 15세이상 인구'
$ws.Range("H274").Value = '조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자'
$ws.Range("I274").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J274").Value = '평균 : 100만원, 중위값 : 200만원'
$ws.Range("K274").Value = '"19.1%"'
$ws.Range("L274").Value = 'Black'
$ws.Range("A275").Value = 45599.521097604171
$ws.Range("B275").Value = 'minheart7844@gmail.com'
$ws.Range("C275").Value = '미디어스쿨'
$ws.Range("D275").Value = 20242503
$ws.Range("E275").Value = '권민'
$ws.Range("F275").Value = '대한민국'
$ws.Range("G275").Value = '취업자 / 경제활동인구'
$ws.Range("H275").Value = '조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자'
$ws.Range("I275").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J275").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K275").Value = '"19.1%"'
$ws.Range("L275").Value = 'Red'
$ws.Range("M275").Value = '나는 사후 장기기증에 참여할 뜻이 없다'
$ws.Range("A276").Value = 45599.534281782406
$ws.Range("B276").Value = 'ehddn0504@naver.com'
$ws.Range("C276").Value = '소프트웨어학부'
$ws.Range("D276").Value = 20245186
$ws.Range("E276").Value = '신동우'
$ws.Range("F276").Value = '스페인'
$ws.Range("G276").Value = '취업자 / 15세 이상 인구'
$ws.Range("H276").Value = '조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자'
$ws.Range("I276").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J276").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K276").Value = '"19.1%"'
$ws.Range("L276").Value = 'Black'
$ws.Range("N276").Value = '나는 사후 장기기증에 참여할 뜻이 있다'
$ws.Range("A277").Value = 45599.543915219911
$ws.Range("B277").Value = 'wlgus4770752@naver.com'
$ws.Range("C277").Value = '의예과'
$ws.Range("D277").Value = 20236121
$ws.Range("E277").Value = '김지현'
$ws.Range("F277").Value = '대한민국'
$ws.Range("G277").Value = '취업자 / 15세 이상 인구'
$ws.Range("H277").Value = '조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자'
$ws.Range("I277").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("J277").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K277").Value = '"19.1%"'
$ws.Range("L277").Value = 'Red'
$ws.Range("M277").Value = '나는 사후 장기기증에 참여할 뜻이 없다'
$ws.Range("A278").Value = 45599.546548229162
$ws.Range("B278").Value = 'kjinju0518@naver.com'
$ws.Range("C278").Value = '경영학과'
$ws.Range("D278").Value = 20232937
$ws.Range("E278").Value = '김진주'
$ws.Range("F278").Value = '대한민국'
$ws.Range("G278").Value = '실업자 / 경제활동인구'
$ws.Range("H278").Value = '조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자'
$ws.Range("I278").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J278").Value = '평균 : 100만원, 중위값 : 200만원'
$ws.Range("K278").Value = '"15%"'
$ws.Range("L278").Value = 'Black'
$ws.Range("N278").Value = '나는 사후 장기기증에 참여할 뜻이 있다'
$ws.Range("A279").Value = 45599.549399791664
$ws.Range("B279").Value = 'bagj11532@gmail.com'
$ws.Range("C279").Value = '체육학과'
$ws.Range("D279").Value = 20244120
$ws.Range("E279").Value = '박준형'
$ws.Range("F279").Value = '대한민국'
$ws.Range("G279").Value = '실업자 / 경제활동인구'
$ws.Range("H279").Value = '조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자'
$ws.Range("I279").Value = '평균 : 100만원, 중위값 : 1,000만원'
$ws.Range("J279").Value = '평균 : 100만원, 중위값 : 1,000만원'
$ws.Range("K279").Value = '"25%"'
$ws.Range("L279").Value = 'Red'
$ws.Range("A280").Value = 45599.551716354166
$ws.Range("B280").Value = 'abcchocoo111@gmail.com'
$ws.Range("C280").Value = '데이터사이언스학부'
$ws.Range("D280").Value = 20243206
$ws.Range("E280").Value = '김기원'
$ws.Range("F280").Value = '대한민국'
$ws.Range("G280").Value = '취업자 / 15세 이상 인구'
$ws.Range("H280").Value = '조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자'
$ws.Range("I280").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J280").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("K280").Value = '"19.1%"'
$ws.Range("L280").Value = 'Red'
$ws.Range("M280").Value = '나는 사후 장기기증에 참여할 뜻이 없다'
$ws.Range("A281").Value = 45599.5607603588
$ws.Range("B281").Value = 'ataraxia050508@naver.com'
$ws.Range("C281").Value = '심리학과'
$ws.Range("D281").Value = 20242118
$ws.Range("E281").Value = '박정호'
$ws.Range("F281").Value = '대한민국'
$ws.Range("G281").Value = '경제활동인구 / 15세이상 인구'
$ws.Range("H281").Value = '조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자'
$ws.Range("I281").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J281").Value = '평균 : 100만원, 중위값 : 1,000만원'
$ws.Range("K281").Value = '"15%"'
$ws.Range("L281").Value = 'Black'
$ws.Range("N281").Value = '나는 사후 장기기증에 참여할 뜻이 있다'
$ws.Range("A282").Value = 45599.577123449075
$ws.Range("B282").Value = 'sshee718@gmail.com'
$ws.Range("C282").Value = '환경생명공학과'
$ws.Range("D282").Value = 20243702
$ws.Range("E282").Value = '권도운'
$ws.Range("F282").Value = '대한민국'
$ws.Range("G282").Value = '취업자 / 15세 이상 인구'
$ws.Range("H282").Value = '조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자'
$ws.Range("I282").Value = '평균 : 100만원, 중위값 : 200만원'
$ws.Range("J282").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K282").Value = '"15%"'
$ws.Range("L282").Value = 'Black'
$ws.Range("A283").Value = 45599.582150578703
$ws.Range("B283").Value = 'dongkyo4@gmail.com'
$ws.Range("C283").Value = '데이터테크'
$ws.Range("D283").Value = 20213241
$ws.Range("E283").Value = '이동교'
$ws.Range("F283").Value = '대한민국'
$ws.Range("G283").Value = '취업자 / 15세 이상 인구'
$ws.Range("H283").Value = '조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자'
$ws.Range("I283").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J283").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K283").Value = '"19.1%"'
$ws.Range("L283").Value = 'Red'

# ============================================================
# Step 3: delete cell M271 entirely (it is no longer the last
# row of the table, so the trailing blank-style cell goes away)
# ============================================================
$ws.Range("M271").Clear()

# ============================================================
# Step 4: resize the Form_Responses1 table to cover the newly
# added rows (A1:N271 -> A1:N283)
# ============================================================
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N283"))

# ============================================================
# Step 5: move the active selection the way the author left it
# ============================================================
$ws.Range("F294").Select()
